$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# New values for columns C (Prophet Forecast), D (Amazon Mean Forecast),
# E (Amazon P70 Forecast), F (Amazon P80 Forecast), G (Amazon P90 Forecast)
# for rows 2-17, after removing Auto Arima from the forecast ensemble.
$data = @(
    @(7, 7, 9, 11, 15),
    @(10, 6, 8, 10, 14),
    @(12, 5, 7, 8, 11),
    @(11, 6, 8, 11, 15),
    @(9, 6, 8, 10, 15),
    @(7, 7, 8, 11, 15),
    @(7, 7, 8, 11, 16),
    @(8, 6, 8, 11, 16),
    @(9, 6, 8, 11, 15),
    @(10, 6, 7, 10, 15),
    @(11, 6, 8, 11, 16),
    @(13, 7, 9, 12, 18),
    @(13, 7, 8, 11, 17),
    @(10, 6, 7, 11, 16),
    @(8, 7, 8, 11, 17),
    @(7, 6, 7, 10, 16)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}

# Update Summary sheet totals derived from the Prophet Forecast column.
# These cells hold numeric-looking text (not real numbers) in the workbook,
# so force a text number format before assigning the value.
$summaryCells = @("B9", "B10", "B11", "B12", "B14")
$summaryValues = @("152", "71", "40", "13", "7")
for ($i = 0; $i -lt $summaryCells.Count; $i++) {
    $cell = $summary.Range($summaryCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $summaryValues[$i]
}
